# Daily attendance processing - 2026-01-27 16:45:49
# Swap the order of names in the "Recorded By" column (G) wherever the
# cell contains both "dnasr281@gmail.com" and "System" separated by a comma.
# "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
